# "22-11 CCodes com reset"
# Adds a new i_Ccload control signal column (T) to the Control-Unit
# Signals truth table, mirroring the existing i_DLen column's header
# style and the per-row 0/1 values used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column T: "i_Ccload" -------------------------------------------
# Header cell (row 4), matches the rest of the N4:S4 header row.
$ws.Range("T4").Value = "i_Ccload"

# Data values for rows 5-19 (one per instruction).
$ccload = @{
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
}
foreach ($r in $ccload.Keys) {
    $ws.Range("T$r").Value = $ccload[$r]
}

# Match the boxed formatting already used across the rest of the header
# row (S4) for the whole new T4:T19 column in one shot.
$ws.Range("S4").Copy()
$ws.Range("T4:T19").PasteSpecial(-4122)
[void]$ws.Application.CutCopyMode

# --- Minor cosmetic touch-ups recorded in the same commit ---------------
$ws.Rows.Item(4).RowHeight = 13.8

$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.Zoom = 120
[void]$ws.Range("G24").Select()
